{"js": "// Fix the \"Luisa Gargano\" scenario sentence:\n//   - \"registrato\"   -> \"registrata\"   (agreement with \"Luisa\", feminine)\n//   - \"dipartimento\" -> \"Dipartimento\" (capitalize department name)\n//\n// The original run's text is replaced in place so the run's character\n// formatting (color / underline) defined on <w:rPr> is preserved.\n\nconst oldSentence =\n  \"Luisa Gargano accede al suo account, precedentemente registrato, e chiede al sistema di visualizzare tutti i report generati dai docenti del dipartimento di Informatica.\";\nconst newSentence =\n  \"Luisa Gargano accede al suo account, precedentemente registrata, e chiede al sistema di visualizzare tutti i report generati dai docenti del Dipartimento di Informatica.\";\n\nconst body = context.document.body;\nconst results = body.search(oldSentence, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found in document body.\");\n}\n\n// Replace the whole sentence in-place (single match expected); insertText\n// with \"Replace\" on the found range keeps the existing run formatting.\nresults.items[0].insertText(newSentence, \"Replace\");\nawait context.sync();\n", "ps1": "# Fix the \"Luisa Gargano\" scenario sentence:\n#   - \"registrato\"   -> \"registrata\"   (agreement with \"Luisa\", feminine)\n#   - \"dipartimento\" -> \"Dipartimento\" (capitalize department name)\n#\n# Use Find/Replace scoped to the full, unique sentence so only this\n# occurrence is touched (plain \"registrato\" / \"dipartimento\" also appear\n# elsewhere in the document with different meaning/casing).\n\n$d = $word.ActiveDocument\n\n$oldSentence = \"Luisa Gargano accede al suo account, precedentemente registrato, e chiede al sistema di visualizzare tutti i report generati dai docenti del dipartimento di Informatica.\"\n$newSentence = \"Luisa Gargano accede al suo account, precedentemente registrata, e chiede al sistema di visualizzare tutti i report generati dai docenti del Dipartimento di Informatica.\"\n\n$rng = $d.Content\n$found = $rng.Find.Execute($oldSentence, $false, $true, $false, $false, $false, $true, 1, $false, $newSentence, 2)\n\nif (-not $found) {\n    throw \"Target sentence not found in document.\"\n}\n"}
